$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "39.876.35"
$cell.Style = $origStyle
$cell = $ws.Range("E2")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.43%  "
$cell.Style = $origStyle

$cell = $ws.Range("D3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.218.97"
$cell.Style = $origStyle
$cell = $ws.Range("E3")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.18%  "
$cell.Style = $origStyle

$cell = $ws.Range("E4")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.08%  "
$cell.Style = $origStyle

$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "292.48"
$cell.Style = $origStyle
$cell = $ws.Range("E5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.79%  "
$cell.Style = $origStyle

$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "86.97"
$cell.Style = $origStyle
$cell = $ws.Range("E6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +6.62%  "
$cell.Style = $origStyle

$cell = $ws.Range("E7")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.94%  "
$cell.Style = $origStyle

$cell = $ws.Range("E8")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.06%  "
$cell.Style = $origStyle

$cell = $ws.Range("E9")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.34%  "
$cell.Style = $origStyle

$cell = $ws.Range("D10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "30.44"
$cell.Style = $origStyle
$cell = $ws.Range("E10")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.83%  "
$cell.Style = $origStyle

$cell = $ws.Range("D11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0785"
$cell.Style = $origStyle
$cell = $ws.Range("E11")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.26%  "
$cell.Style = $origStyle

$cell = $ws.Range("D12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "47.60"
$cell.Style = $origStyle
$cell = $ws.Range("E12")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.09%  "
$cell.Style = $origStyle

$cell = $ws.Range("E13")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.73%  "
$cell.Style = $origStyle

$cell = $ws.Range("D14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "6.34"
$cell.Style = $origStyle
$cell = $ws.Range("E14")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.03%  "
$cell.Style = $origStyle

$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.561.84"
$cell.Style = $origStyle
$cell = $ws.Range("E15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.92%  "
$cell.Style = $origStyle

$cell = $ws.Range("D16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "14.02"
$cell.Style = $origStyle
$cell = $ws.Range("E16")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +0.66%  "
$cell.Style = $origStyle

$cell = $ws.Range("D17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.228.60"
$cell.Style = $origStyle
$cell = $ws.Range("E17")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.44%  "
$cell.Style = $origStyle

$cell = $ws.Range("D18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.730"
$cell.Style = $origStyle
$cell = $ws.Range("E18")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.02%  "
$cell.Style = $origStyle

$cell = $ws.Range("D19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "39.812.06"
$cell.Style = $origStyle
$cell = $ws.Range("E19")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.61%  "
$cell.Style = $origStyle

$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0881"
$cell.Style = $origStyle
$cell = $ws.Range("E20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.64%  "
$cell.Style = $origStyle

$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "11.27"
$cell.Style = $origStyle
$cell = $ws.Range("E21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +9.97%  "
$cell.Style = $origStyle

$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "5.80"
$cell.Style = $origStyle
$cell = $ws.Range("E22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.66%  "
$cell.Style = $origStyle

$cell = $ws.Range("D23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "65.66"
$cell.Style = $origStyle
$cell = $ws.Range("E23")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.54%  "
$cell.Style = $origStyle

$cell = $ws.Range("D24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "235.82"
$cell.Style = $origStyle
$cell = $ws.Range("E24")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.00%  "
$cell.Style = $origStyle

$cell = $ws.Range("E25")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.30%  "
$cell.Style = $origStyle

$cell = $ws.Range("E26")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.10%  "
$cell.Style = $origStyle

$cell = $ws.Range("E27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.26%  "
$cell.Style = $origStyle

$cell = $ws.Range("D28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "22.77"
$cell.Style = $origStyle
$cell = $ws.Range("E28")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.17%  "
$cell.Style = $origStyle

$cell = $ws.Range("E29")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.15%  "
$cell.Style = $origStyle

$cell = $ws.Range("D30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "9.27"
$cell.Style = $origStyle
$cell = $ws.Range("E30")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.37%  "
$cell.Style = $origStyle

$cell = $ws.Range("D31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "32.84"
$cell.Style = $origStyle
$cell = $ws.Range("E31")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.80%  "
$cell.Style = $origStyle

$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "152.02"
$cell.Style = $origStyle
$cell = $ws.Range("E32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.37%  "
$cell.Style = $origStyle

$cell = $ws.Range("E33")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.11%  "
$cell.Style = $origStyle

$cell = $ws.Range("D34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "4.94"
$cell.Style = $origStyle
$cell = $ws.Range("E34")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.37%  "
$cell.Style = $origStyle

$cell = $ws.Range("E35")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.05%  "
$cell.Style = $origStyle

$cell = $ws.Range("D36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.38"
$cell.Style = $origStyle
$cell = $ws.Range("E36")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.19%  "
$cell.Style = $origStyle

$cell = $ws.Range("B37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "Stellar"
$cell.Style = $origStyle
$cell = $ws.Range("C37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$cell.Style = $origStyle
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.112"
$cell.Style = $origStyle
$cell = $ws.Range("E37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.16%  "
$cell.Style = $origStyle

$cell = $ws.Range("B38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "LidoDAOToken"
$cell.Style = $origStyle
$cell = $ws.Range("C38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$cell.Style = $origStyle
$cell = $ws.Range("D38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.79"
$cell.Style = $origStyle
$cell = $ws.Range("E38")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +6.57%  "
$cell.Style = $origStyle

$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "15.89"
$cell.Style = $origStyle
$cell = $ws.Range("E39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.77%  "
$cell.Style = $origStyle

$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "0.0991"
$cell.Style = $origStyle
$cell = $ws.Range("E40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.61%  "
$cell.Style = $origStyle

$cell = $ws.Range("E41")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +4.16%  "
$cell.Style = $origStyle

$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "3.79"
$cell.Style = $origStyle
$cell = $ws.Range("E42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +5.51%  "
$cell.Style = $origStyle

$cell = $ws.Range("D43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.062.74"
$cell.Style = $origStyle
$cell = $ws.Range("E43")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +8.81%  "
$cell.Style = $origStyle

$cell = $ws.Range("E44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +11.95%  "
$cell.Style = $origStyle

$cell = $ws.Range("E45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +3.84%  "
$cell.Style = $origStyle

$cell = $ws.Range("E46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +11.77%  "
$cell.Style = $origStyle

$cell = $ws.Range("E47")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.16%  "
$cell.Style = $origStyle

$cell = $ws.Range("D48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.62"
$cell.Style = $origStyle
$cell = $ws.Range("E48")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.38%  "
$cell.Style = $origStyle

$cell = $ws.Range("D49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "2.435.67"
$cell.Style = $origStyle
$cell = $ws.Range("E49")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +1.20%  "
$cell.Style = $origStyle

$cell = $ws.Range("D50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "71.17"
$cell.Style = $origStyle
$cell = $ws.Range("E50")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  -0.27%  "
$cell.Style = $origStyle

$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "89.21"
$cell.Style = $origStyle
$cell = $ws.Range("E51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = "  +2.61%  "
$cell.Style = $origStyle
